# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns,
# and fix the MXToken/FraxShare row order (rows 49-50).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.727.59"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.105.61"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.03"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.18"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.393"
$ws.Range("E9").Value = "  +2.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0778"
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.413.58"
$ws.Range("E12").Value = "  +1.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.57"
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.14"
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.23"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.092.66"
$ws.Range("E17").Value = "  +1.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.662.78"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.22"
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.04"
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.47"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.40"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("E26").Value = "  +1.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.97"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("E28").Value = "  +2.82%  "
$ws.Range("E29").Value = "  -4.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.47"
$ws.Range("E30").Value = "  +1.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.118"
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.64"
$ws.Range("E32").Value = "  +3.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0622"
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.57"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.59"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.49"
$ws.Range("E36").Value = "  +5.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.77"
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.61"
$ws.Range("E39").Value = "  -6.63%  "
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0964"
$ws.Range("E41").Value = "  +1.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.45"
$ws.Range("E42").Value = "  +2.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.472.61"
$ws.Range("E43").Value = "  +0.95%  "
$ws.Range("E44").Value = "  +1.02%  "
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.16"
$ws.Range("E46").Value = "  -11.44%  "
$ws.Range("E47").Value = "  +1.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.59"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.32"
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.03"
$ws.Range("E50").Value = "  +3.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.299.57"
$ws.Range("E51").Value = "  +1.68%  "
